# Apply the two changes captured by the commit:
#  1. Slide 16's table switches from the custom "Table_0" table style
#     to the built-in PowerPoint table style {F93B38A1-EDB9-4F2F-9723-6F73325CA2C2}.
#  2. The (slide) master's theme colour scheme is changed from the
#     "Integral" palette to the default "Office" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on Slide 16, Shape 3 (the graphicFrame/table) ---
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{F93B38A1-EDB9-4F2F-9723-6F73325CA2C2}")

# --- 2. Theme colour scheme: Integral -> Office (default) colours ---
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0          # Dark 1    -> 000000
$colorScheme.Item(2).RGB  = 16777215   # Light 1   -> FFFFFF
$colorScheme.Item(3).RGB  = 6968388    # Dark 2    -> 44546A
$colorScheme.Item(4).RGB  = 15132391   # Light 2   -> E7E6E6
$colorScheme.Item(5).RGB  = 13998939   # Accent 1  -> 5B9BD5
$colorScheme.Item(6).RGB  = 3243501    # Accent 2  -> ED7D31
$colorScheme.Item(7).RGB  = 10855845   # Accent 3  -> A5A5A5
$colorScheme.Item(8).RGB  = 49407      # Accent 4  -> FFC000
$colorScheme.Item(9).RGB  = 12874308   # Accent 5  -> 4472C4
$colorScheme.Item(10).RGB = 4697456    # Accent 6  -> 70AD47
$colorScheme.Item(11).RGB = 12673797   # Hyperlink -> 0563C1
$colorScheme.Item(12).RGB = 7491477    # Followed Hyperlink -> 954F72
